$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: student "anne" -> "nosha", and her numeric id column cleaned up
$ws.Range("A3").Value = "nosha"
$ws.Range("E3").Value = "'22102"
$ws.Range("E3").ClearFormats()

# Row 4: "dodo" row - just the id column format cleanup (22101.0 -> 22101)
$ws.Range("E4").Value = "'22101"
$ws.Range("E4").ClearFormats()

# Row 5 ("lola") is removed entirely - delete the whole row, shifting rows below up
$ws.Rows("5").Delete()

# The row that used to be 6 ("3li") is now row 5 - clean up its id column too
$ws.Range("E5").Value = "'22103"
$ws.Range("E5").ClearFormats()
